$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Workbook / sheet level
# ---------------------------------------------------------------------------
$ws.Name = "Putz 101"

# ---------------------------------------------------------------------------
# Column widths (character units). The engine quantises ColumnWidth to the
# nearest 1/6 of a character, so these are the closest achievable values to
# the authored widths 5.44140625 / 10.33203125 / 12.77734375 / 3.21875 /
# 27.44140625.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 4.666666666666667
$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 2.333333333333333
$ws.Columns.Item(5).ColumnWidth = 26.66666666666667

# ---------------------------------------------------------------------------
# Values
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "DC  no:23-24QPutz101"
$ws.Range("E1").Value = "date: 17/04/2023"

$ws.Range("A2").Value = "client: PUTZMEISTER CONCRETE MACHINES PVT LTD"
$ws.Range("E2").Value = "Ref: Replacement towards invoice number - b23-24MQ101"

$ws.Range("A3").Value = "Billing Address: Verna Goa"

$ws.Range("E4").Value = "gstin: 30BBPPP8605M1Z7"

$ws.Range("A5").Value = "Shipping Address: Verna Goa"

$ws.Range("A7").Value = "SrNo."
$ws.Range("B7").Value = "Part No."
$ws.Range("C7").Value = "Item Description"
$ws.Range("D7").Value = "Qty"

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 616026
$ws.Range("C8").Value = "Screw Conveyor Bellow_Dia 224x500mm"
$ws.Range("D8").Value = 6

$ws.Range("A10").Value = "Terms and conditions apply as per norms JPT001"

# ---------------------------------------------------------------------------
# Merges
# ---------------------------------------------------------------------------
$ws.Range("A1:C1").Merge()
$ws.Range("A2:C2").Merge()
$ws.Range("A3:C3").Merge()
$ws.Range("A5:C5").Merge()

# ---------------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 20.4
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 35.4
$ws.Rows.Item(5).RowHeight = 37.2
$ws.Rows.Item(7).RowHeight = 28.8
$ws.Rows.Item(8).RowHeight = 61.8
$ws.Rows.Item(10).RowHeight = 15

# ---------------------------------------------------------------------------
# Row 1: A1:C1 box (left,top,bottom) | D1 top only | E1 full box
# ---------------------------------------------------------------------------
$ws.Range("A1:C1").Borders.Item(7).LineStyle = 1
$ws.Range("A1:C1").Borders.Item(8).LineStyle = 1
$ws.Range("A1:C1").Borders.Item(9).LineStyle = 1
$ws.Range("A1:C1").HorizontalAlignment = -4131
$ws.Range("A1:C1").VerticalAlignment = -4160

$ws.Range("D1").Borders.Item(8).LineStyle = 1

$ws.Range("E1").Borders.Item(7).LineStyle = 1
$ws.Range("E1").Borders.Item(8).LineStyle = 1
$ws.Range("E1").Borders.Item(9).LineStyle = 1
$ws.Range("E1").Borders.Item(10).LineStyle = 1
$ws.Range("E1").HorizontalAlignment = -4131
$ws.Range("E1").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Row 2: A2:C2 box (top,bottom only) | D2 none | E2 full box
# ---------------------------------------------------------------------------
$ws.Range("A2:C2").Borders.Item(8).LineStyle = 1
$ws.Range("A2:C2").Borders.Item(9).LineStyle = 1
$ws.Range("A2:C2").HorizontalAlignment = -4131
$ws.Range("A2:C2").VerticalAlignment = -4108
$ws.Range("A2:C2").WrapText = $true

$ws.Range("E2").Borders.Item(7).LineStyle = 1
$ws.Range("E2").Borders.Item(8).LineStyle = 1
$ws.Range("E2").Borders.Item(9).LineStyle = 1
$ws.Range("E2").Borders.Item(10).LineStyle = 1
$ws.Range("E2").HorizontalAlignment = -4131
$ws.Range("E2").VerticalAlignment = -4108
$ws.Range("E2").WrapText = $true

# ---------------------------------------------------------------------------
# Row 3: A3:C3 box (left,top,bottom) | D3 none | E3 left only
# ---------------------------------------------------------------------------
$ws.Range("A3:C3").Borders.Item(7).LineStyle = 1
$ws.Range("A3:C3").Borders.Item(8).LineStyle = 1
$ws.Range("A3:C3").Borders.Item(9).LineStyle = 1
$ws.Range("A3:C3").HorizontalAlignment = -4131
$ws.Range("A3:C3").VerticalAlignment = -4160
$ws.Range("A3:C3").WrapText = $true

$ws.Range("E3").Borders.Item(7).LineStyle = 1

# ---------------------------------------------------------------------------
# Row 4: A4 left only | E4 full box
# ---------------------------------------------------------------------------
$ws.Range("A4").Borders.Item(7).LineStyle = 1

$ws.Range("E4").Borders.Item(7).LineStyle = 1
$ws.Range("E4").Borders.Item(8).LineStyle = 1
$ws.Range("E4").Borders.Item(9).LineStyle = 1
$ws.Range("E4").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------------
# Row 5: A5:C5 box (left,top,bottom) | D5 bottom only | E5 right+bottom
# ---------------------------------------------------------------------------
$ws.Range("A5:C5").Borders.Item(7).LineStyle = 1
$ws.Range("A5:C5").Borders.Item(8).LineStyle = 1
$ws.Range("A5:C5").Borders.Item(9).LineStyle = 1
$ws.Range("A5:C5").HorizontalAlignment = -4131
$ws.Range("A5:C5").VerticalAlignment = -4160
$ws.Range("A5:C5").WrapText = $true

$ws.Range("D5").Borders.Item(9).LineStyle = 1

$ws.Range("E5").Borders.Item(9).LineStyle = 1
$ws.Range("E5").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------------
# Row 7: header row - A7:D7 each a full box, bold, centered, wrap
# ---------------------------------------------------------------------------
$ws.Range("A7:D7").Borders.Item(7).LineStyle = 1
$ws.Range("A7:D7").Borders.Item(8).LineStyle = 1
$ws.Range("A7:D7").Borders.Item(9).LineStyle = 1
$ws.Range("A7:D7").Borders.Item(10).LineStyle = 1
$ws.Range("A7:D7").HorizontalAlignment = -4108
$ws.Range("A7:D7").VerticalAlignment = -4108
$ws.Range("A7:D7").WrapText = $true
$ws.Range("A7:D7").Font.Bold = $true

# ---------------------------------------------------------------------------
# Row 8: data row - A8:D8 each a full box, centered, wrap
# ---------------------------------------------------------------------------
$ws.Range("A8:D8").Borders.Item(7).LineStyle = 1
$ws.Range("A8:D8").Borders.Item(8).LineStyle = 1
$ws.Range("A8:D8").Borders.Item(9).LineStyle = 1
$ws.Range("A8:D8").Borders.Item(10).LineStyle = 1
$ws.Range("A8:D8").HorizontalAlignment = -4108
$ws.Range("A8:D8").VerticalAlignment = -4108
$ws.Range("A8:D8").WrapText = $true

$ws.Range("B8:C8").Font.Name = "Arial"
$ws.Range("B8:C8").Font.Size = 10

# ---------------------------------------------------------------------------
# Selection (cosmetic, matches the authored cursor position)
# ---------------------------------------------------------------------------
$ws.Range("I8").Select()
